$d = $word.ActiveDocument
$wdParagraph = 4

# ---------------------------------------------------------------------------
# 1) Remove the three "Responsáveis pela Elaboração" / author-name runs, but
#    keep their (now empty) paragraphs.
# ---------------------------------------------------------------------------

# "Responsáveis pela Elaboração" -> just delete the run's text
$rng = $d.Content
if ($rng.Find.Execute("Responsáveis pela Elaboração")) {
    $rng.Delete()
}

# "Moisés Hilário Rodrigues" -> just delete the run's text
$rng = $d.Content
if ($rng.Find.Execute("Moisés Hilário Rodrigues")) {
    $rng.Delete()
}

# "Igor Moura Brandão" -> the paragraph mark itself also carries the bold /
# green / sz24 formatting in its pPr/rPr, so replace the whole paragraph
# (mark included) with a clean one that only keeps the spacing/jc/rFonts
# that are shared with the sibling paragraphs.
$rng = $d.Content
if ($rng.Find.Execute("Igor Moura Brandão")) {
    $rng.Expand($wdParagraph)
    $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:spacing w:before='60' w:after='60'/><w:jc w:val='center'/><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/></w:rPr></w:pPr></w:p>"
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 2) Insert a brand-new empty paragraph carrying the "_GoBack" bookmark right
#    after the (now empty) paragraph that used to follow "Igor Moura
#    Brandão" - i.e. immediately before the "Detalhes do Caso de Teste"
#    heading.
# ---------------------------------------------------------------------------

$rng = $d.Content
if ($rng.Find.Execute("Detalhes do C")) {
    $rng.Expand($wdParagraph)
    $headingStart = $rng.Start
    $precedingMark = $d.Range($headingStart - 1, $headingStart)
    $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/></w:rPr></w:pPr></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/></w:rPr></w:pPr><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
    $precedingMark.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3) The trailing paragraph (right after the last table, before the final
#    sectPr) no longer owns the "_GoBack" bookmark - it becomes a plain
#    empty paragraph.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $lastPara.Range
$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>"
$lastRange.InsertXML($xml)

Write-Output "done"
